$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 62 (shifts existing rows 62:189 down to 63:190
# and grows the used range from A1:R189 to A1:R190).
$ws.Rows(62).Insert()

# Populate the newly inserted row 62 with a new daily price record.
$ws.Range("A62").Value = 4
$ws.Range("B62").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C62").Value = "Los Lagos"
$ws.Range("D62").Value = 44536
$ws.Range("E62").Value = 10
$ws.Range("F62").Value = 100112043
$ws.Range("G62").Value = "Pepino ensalada"
$ws.Range("H62").Value = "Sin especificar"
$ws.Range("I62").Value = "Primera"
$ws.Range("J62").Value = 200
$ws.Range("K62").Value = 11000
$ws.Range("L62").Value = 12000
$ws.Range("M62").Value = 11500
$ws.Range("N62").Value = "`$/caja 60 unidades"
$ws.Range("O62").Value = "Región de Arica y Parinacota"
$ws.Range("P62").Value = 192
$ws.Range("Q62").Value = 60
$ws.Range("R62").Value = "Hortaliza"
